# DAY_01 log was reset: the rows that tracked the "Tomcat / Maven / README /
# Git" tasks (now tracked elsewhere) are cleared back to blank template rows,
# which also drops the now-orphaned shared strings and the hyperlink that
# lived in E2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAY_01")

# Clear out the S.No / Date / Task / links / time / error columns for the
# four data rows (2-5). Row 5's D:E cells were already empty in the source.
$ws.Range("A2:G4").ClearContents()
$ws.Range("A5:C5").ClearContents()

# E2 carried an external hyperlink to the eclipse.org Tomcat tutorial; once
# the cell text is gone the link itself needs to go too.
$ws.Hyperlinks.Delete()

# Leave the sheet active with the cursor back at the top of the table
# instead of scrolled over to F5.
$ws.Activate()
$ws.Range("A3").Select()
